$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the data block (rows 2 and 3), pushing
# all the existing weekly records down by two rows.
$ws.Rows("2:3").Insert()

# The inserted rows pick up inherited formatting from the row above (the
# bold header). Clear that back to the default, then restore the date
# number format on column D to match the rest of the table.
$ws.Range("A2:T3").ClearFormats()
$ws.Range("D2:D3").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Columns that are constant for every record in this sheet.
$ws.Range("A2:A3").Value = 11
$ws.Range("B2:B3").Value = "Vega Monumental Concepción"
$ws.Range("C2:C3").Value = "Bíobío"
$ws.Range("E2:E3").Value = 8
$ws.Range("F2:F3").Value = "Fruta"
$ws.Range("G2:G3").Value = 100107
$ws.Range("H2:H3").Value = "Otros"
$ws.Range("I2:I3").Value = 100107011
$ws.Range("J2:J3").Value = "Tuna"
$ws.Range("K2:K3").Value = "Sin especificar"
$ws.Range("R2:R3").Value = "Provincia de Melipilla"
$ws.Range("T2:T3").Value = 18

# New row 2: "Primera" quality record for the latest week.
$ws.Range("D2").Value = 45168
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 50
$ws.Range("N2").Value = 26000
$ws.Range("O2").Value = 26000
$ws.Range("P2").Value = 26000
$ws.Range("Q2").Value = "$/caja 18 kilos"
$ws.Range("S2").Value = 1444

# New row 3: "Segunda" quality record for the latest week.
$ws.Range("D3").Value = 45168
$ws.Range("L3").Value = "Segunda"
$ws.Range("M3").Value = 50
$ws.Range("N3").Value = 22000
$ws.Range("O3").Value = 22000
$ws.Range("P3").Value = 22000
$ws.Range("Q3").Value = "$/caja 18 kilos"
$ws.Range("S3").Value = 1222

Write-Output $ws.UsedRange.Address()
